$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("I2").Value = 7
$ws.Range("J2").Value = 9

$ws.Range("I3").Value = 5
$ws.Range("J3").Value = 7

$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 4

$ws.Range("I5").Value = 1
$ws.Range("J5").Value = 4
